# Update cryptos list with latest scraped prices/volumes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.647.60"
$ws.Range("E2").Value = "  -0.43%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.849.49"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -1.57%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.03"
$ws.Range("E5").Value = "  -1.33%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -1.30%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4225"
$ws.Range("E7").Value = "  +0.27%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3640"
$ws.Range("E8").Value = "  +0.49%  "

# Row 9 - OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.45"
$ws.Range("E9").Value = "  -0.83%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07285"
$ws.Range("E10").Value = "  +1.75%  "

# Row 11 - Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8743"
$ws.Range("E11").Value = "  -2.44%  "

# Row 12 - Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.65"
$ws.Range("E12").Value = "  +0.93%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.817.81"
$ws.Range("E13").Value = "  -2.17%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.336"
$ws.Range("E14").Value = "  +1.09%  "

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.515"
$ws.Range("E15").Value = "  -0.27%  "

# Row 16 - TRON
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06854"
$ws.Range("E16").Value = "  +0.00%  "

# Row 17 - BinanceUSD
$ws.Range("E17").Value = "  -1.34%  "

# Row 18 - Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.46"
$ws.Range("E18").Value = "  +2.37%  "

# Row 19 - ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008926"
$ws.Range("E19").Value = "  +3.10%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  -0.91%  "

# Row 21 - Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.33"
$ws.Range("E21").Value = "  +0.21%  "

# Row 22 - WrappedBTC
$ws.Range("D22").Value = "27.672.95"
$ws.Range("E22").Value = "  -0.14%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.986"
$ws.Range("E23").Value = "  +1.05%  "

# Row 24 - Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.35"
$ws.Range("E24").Value = "  -4.04%  "

# Row 25 - WrappedliquidstakedEther2.0
$ws.Range("D25").Value = "2.087.51"
$ws.Range("E25").Value = "  +0.06%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.84"
$ws.Range("E27").Value = "  +0.99%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +5.66%  "

# Row 29 - BitcoinCash
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.93"
$ws.Range("E29").Value = "  +10.86%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.248"
$ws.Range("E30").Value = "  -0.56%  "

# Row 31 - LidoDAOToken
$ws.Range("E31").Value = "  +15.85%  "

# Row 32 - Stellar
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08854"
$ws.Range("E32").Value = "  -0.50%  "

# Row 33 - ImmutableX
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7664"
$ws.Range("E33").Value = "  -0.26%  "

# Row 34 - Filecoin->HuobiToken
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.976"
$ws.Range("E34").Value = "  +3.55%  "

# Row 35 - HuobiToken->Filecoin
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.538"
$ws.Range("E35").Value = "  +1.27%  "

# Row 36 - ARBITRUM
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.103"
$ws.Range("E36").Value = "  +4.20%  "

# Row 37 - Frax
$ws.Range("E37").Value = "  -1.61%  "

# Row 38 - TrustWalletToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.095"
$ws.Range("E38").Value = "  -1.15%  "

# Row 39 - Hedera
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05347"
$ws.Range("E39").Value = "  +0.80%  "

# Row 40 - VeChain
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01928"
$ws.Range("E40").Value = "  +1.01%  "

# Row 41 - MXToken
$ws.Range("E41").Value = "  -5.04%  "

# Row 42 - TheSandbox->FraxShare
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.887"
$ws.Range("E42").Value = "  +0.59%  "

# Row 43 - FraxShare->TheSandbox
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5086"
$ws.Range("E43").Value = "  +0.52%  "

# Row 44 - Algorand
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1646"
$ws.Range("E44").Value = "  +0.97%  "

# Row 45 - Aptos
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.315"
$ws.Range("E45").Value = "  +0.65%  "

# Row 46 - Cronos
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06534"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47 - Quant->EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("E47").Value = "  +1.45%  "

# Row 48 - EnergySwap->Quant
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.58"
$ws.Range("E48").Value = "  +1.18%  "

# Row 49 - Decentraland
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4693"
$ws.Range("E49").Value = "  +0.01%  "

# Row 50 - PaxDollar
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("E50").Value = "  -1.35%  "

# Row 51 - NEARProtocol
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.623"
$ws.Range("E51").Value = "  +0.80%  "
